$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.052.24"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "2.602.81"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'590.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").Value = "'149.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "2.601.24"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'0.127"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").Value = "'27.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "3.070.64"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.05%  "
$ws.Range("D17").Value = "66.927.12"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "2.602.11"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'363.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "'7.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.74%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "'4.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "'72.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'9.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.733.17"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "'584.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "0.0₃0985"
$ws.Range("E31").Value = "  -6.69%  "
$ws.Range("D32").Value = "'1.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.65%  "
$ws.Range("D33").Value = "'7.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.82%  "
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("D38").Value = "'155.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "'18.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "'5.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("D44").Value = "'17.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'153.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").Value = "0.0₆0283"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").Value = "'21.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.51%  "
